$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: continuation entry for Stroh et. al. (2012) / GC ADPs
$ws.Range("A21").Value = "Stroh et. al. (2012)"
$ws.Range("B21").Value = "NMDA Receptor-Dependent Synaptic Activation of TRPC`nChannels in Olfactory Bulb Granule Cells"
$ws.Range("C21").Value = "Mice"
$ws.Range("D21").Value = "GC"
$ws.Range("E21").Value = "GCs show ADPs in vivo to somatic stimulation, and LLDs to glomerular stimulation"
$ws.Range("B21").Style = $ws.Range("B18").Style

# Row 22: ADP amplitude data row
$ws.Range("A22").Value = "Stroh et. al. (2012)"
$ws.Range("B22").Value = "NMDA Receptor-Dependent Synaptic Activation of TRPC`nChannels in Olfactory Bulb Granule Cells"
$ws.Range("C22").Value = "Mice"
$ws.Range("D22").Value = "GC"
$ws.Range("E22").Value = "ADP amplitude"
$ws.Range("F22").Value = "11.1+-4.7 mV"
$ws.Range("G22").Value = 49
$ws.Range("H22").Value = 21
$ws.Range("I22").Value = "not REPORTED"
$ws.Range("J22").Value = "Soma injections 1000pa for 1ms. a mean sAP-ADP amplitude above resting potential"
$ws.Range("K22").Value = "STD"
$ws.Range("B22").Style = $ws.Range("B18").Style

# Row 23: ADP half-duration data row
$ws.Range("E23").Value = "ADP half-duration"
$ws.Range("F23").Value = "42+-22 ms"
$ws.Range("G23").Value = 49
$ws.Range("H23").Value = 21
$ws.Range("I23").Value = "not REPORTED"
$ws.Range("J23").Value = "Soma injections 1000pa for 1ms"
$ws.Range("K23").Value = "Throughout the paper, tau1/2 denotes halfdurations from the peak amplitude of Vm onward, measured between the onset of the afterdepolarization (ADP) right after the sodium spike and one-half of its maximum amplitude"
$ws.Range("K23").Style = $ws.Range("B18").Style

# The multi-line title text triggers Excel's auto row-height; re-fit so the
# rows stay at the default height (matching row 18's existing behavior).
$ws.Rows.Item(21).AutoFit()
$ws.Rows.Item(22).AutoFit()

# Update the active selection to reflect the new last-used row (mirrors the
# cursor landing one row below the newly added block, same as before).
$ws.Range("E20").Select()
